$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns A and B: unitTypeId moves to column A, projectId moves to column B
$ws.Range("A1").Value = "unitTypeId"
$ws.Range("B1").Value = "projectId"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0

# Update the active selection on the sheet
$ws.Range("B4").Select() | Out-Null
